# Timesheet update: add entries for Feb 28 2020.
#   - Row 230: blank separator row (same look as the row-219 separator
#     that precedes the existing "Feb 27" block).
#   - Rows 231-235: the new Feb 28 timesheet entries, styled the same way
#     as the existing data rows (e.g. row 220 which starts the Feb 27
#     block: column A/C centered, column B left-aligned).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values ---

# Row 230 left blank (it's just the colored separator bar between days).
$ws.Range("A230").Value = ""
$ws.Range("B230").Value = ""
$ws.Range("C230").Value = ""

$ws.Range("A231").Value = "Feb 28 10:00 to 11:00"
$ws.Range("B231").Value = "Stated model tuning, creating model tuning class, writing model tuning parameters."
$ws.Range("C231").Value = "Infimetrics"

$ws.Range("A232").Value = "Feb 28 11:00 to 12:30"
$ws.Range("B232").Value = "Client call"
$ws.Range("C232").Value = "Sapphire automation"

$ws.Range("A233").Value = "Feb 28 12:30 to 13:30"
$ws.Range("B233").Value = "Office anouncement and some documentations. "
$ws.Range("C233").Value = "Infimetrics"

$ws.Range("A234").Value = "Feb 28 13:30 to 14:30"
$ws.Range("B234").Value = "Lunch"
$ws.Range("C234").Value = "Infimetrics"

$ws.Range("A235").Value = "Feb 28 14:30 to 19:00"
$ws.Range("B235").Value = "Python class and office work."
$ws.Range("C235").Value = "Infimetrics"

# --- Formatting ---
# Reuse the existing cell styles (fill colour / alignment) by copying the
# format from the previous day's equivalent rows, rather than building new
# style entries from scratch.

# Separator row: copy from row 219 (the blank separator before "Feb 27").
$ws.Range("A219:C219").Copy()
$ws.Range("A230:C230").PasteSpecial(-4122)

# Data rows: copy from row 220 (first data row of the "Feb 27" block).
$ws.Range("A220:C220").Copy()
$ws.Range("A231:C231").PasteSpecial(-4122)
$ws.Range("A220:C220").Copy()
$ws.Range("A232:C232").PasteSpecial(-4122)
$ws.Range("A220:C220").Copy()
$ws.Range("A233:C233").PasteSpecial(-4122)
$ws.Range("A220:C220").Copy()
$ws.Range("A234:C234").PasteSpecial(-4122)
$ws.Range("A220:C220").Copy()
$ws.Range("A235:C235").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the visible selection down to the newly added last row, matching
# where the author ended up after typing the new entries.
$ws.Range("D235").Select()

# Scroll the window so row 214 is near the top, matching the author's
# final viewport (topLeftCell="A214").
$excel.ActiveWindow.ScrollRow = 214
$excel.ActiveWindow.ScrollColumn = 1
